$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new 2022 column (K) of data
$ws.Range("K4").Value = 2022
$ws.Range("K5").Value = 26.495524312074597
$ws.Range("K6").Value = 59.383769502755833
$ws.Range("K7").Value = 38.32334404557426
$ws.Range("K8").Value = 48.136790950525594
$ws.Range("K9").Value = 46.63213064070051
$ws.Range("K10").Value = 32.657429481680126
$ws.Range("K11").Value = 31.457245964894081
$ws.Range("K12").Value = 22.734405597714229
$ws.Range("K13").Value = -0.19691879995369213
$ws.Range("K14").Value = 33.158040409631916

# Copy styles from column J so K matches existing formatting for each row
$ws.Range("J4:J14").Copy()
$ws.Range("K4").PasteSpecial(-4122)  # xlPasteFormats

$excel.CutCopyMode = $false

# Update the selected cell to reflect the new selection in the diff
$ws.Range("M7").Select()
